$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1242.7142
$ws.Range("I28").Value = 1126.909
$ws.Range("J28").Value = 1667.3334
$ws.Range("K28").Value = 1126.909
$ws.Range("L28").Value = 1667.3334
$ws.Range("M28").Value = -641.9090000000001
$ws.Range("N28").Value = -2637.3334

$ws.Range("H32").Value = 549.75
$ws.Range("J32").Value = 499
$ws.Range("L32").Value = 499
$ws.Range("N32").Value = -1151

$ws.Range("H62").Value = 8502.5
$ws.Range("I62").Value = 7002
$ws.Range("J62").Value = 10003
$ws.Range("K62").Value = 7002
$ws.Range("L62").Value = 10003
$ws.Range("M62").Value = -6378
$ws.Range("N62").Value = -11251

$ws.Range("H65").Value = 8502.5
$ws.Range("I65").Value = 7002
$ws.Range("J65").Value = 10003
$ws.Range("K65").Value = 35010
$ws.Range("L65").Value = 50015
$ws.Range("M65").Value = -31890
$ws.Range("N65").Value = -56255

$ws.Range("H74").Value = 11116534
$ws.Range("I74").Value = 25003500
$ws.Range("J74").Value = 6961.6
$ws.Range("K74").Value = 25003500
$ws.Range("L74").Value = 6961.6
$ws.Range("M74").Value = -25002564
$ws.Range("N74").Value = -8833.6

$ws.Range("H77").Value = 11116534
$ws.Range("I77").Value = 25003500
$ws.Range("J77").Value = 6961.6
$ws.Range("K77").Value = 125017500
$ws.Range("L77").Value = 34808
$ws.Range("M77").Value = -125012820
$ws.Range("N77").Value = -44168

$ws.Range("H138").Value = 2741.37
$ws.Range("I138").Value = 1253.7142
$ws.Range("J138").Value = 2853.344
$ws.Range("K138").Value = 3761.1426
$ws.Range("L138").Value = 8560.032
$ws.Range("M138").Value = 1378.8574
$ws.Range("N138").Value = -18840.032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3623.7124
$ws.Range("I32").Value = 4058.8823
$ws.Range("J32").Value = 2614.9092
$ws.Range("K32").Value = 4058.8823
$ws.Range("L32").Value = 2614.9092
$ws.Range("M32").Value = -3771.8823
$ws.Range("N32").Value = -3188.9092

$ws.Range("H37").Value = 20734.188
$ws.Range("I37").Value = 7350
$ws.Range("J37").Value = 28764.7
$ws.Range("K37").Value = 7350
$ws.Range("L37").Value = 28764.7
$ws.Range("M37").Value = -7077
$ws.Range("N37").Value = -29310.7

$ws.Range("H41").Value = 12408.5
$ws.Range("I41").Value = 2220.6667
$ws.Range("J41").Value = 27690.25
$ws.Range("K41").Value = 2220.6667
$ws.Range("L41").Value = 27690.25
$ws.Range("M41").Value = -1806.6667
$ws.Range("N41").Value = -28518.25

$ws.Range("H45").Value = 2224
$ws.Range("I45").Value = 1958.4
$ws.Range("K45").Value = 1958.4
$ws.Range("M45").Value = -1581.4

$ws.Range("H110").Value = 1284.579
$ws.Range("I110").Value = 1367.3
$ws.Range("J110").Value = 1192.6666
$ws.Range("K110").Value = 1367.3
$ws.Range("L110").Value = 1192.6666
$ws.Range("M110").Value = 677.7
$ws.Range("N110").Value = -5282.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 166667920
$ws.Range("I105").Value = 333333340
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 333333340
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -333331593
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12346947
$ws.Range("I16").Value = 22223344
$ws.Range("J16").Value = 1449.75
$ws.Range("K16").Value = 22223344
$ws.Range("L16").Value = 1449.75
$ws.Range("M16").Value = -22223057
$ws.Range("N16").Value = -2023.75

$ws.Range("H93").Value = 16586.666
$ws.Range("I93").Value = 7203.5
$ws.Range("J93").Value = 35353
$ws.Range("K93").Value = 7203.5
$ws.Range("L93").Value = 35353
$ws.Range("M93").Value = -5331.5
$ws.Range("N93").Value = -39097

$ws.Range("H99").Value = 11115566
$ws.Range("J99").Value = 5483.3335
$ws.Range("L99").Value = 5483.3335
$ws.Range("N99").Value = -8479.3335

$ws.Range("H105").Value = 2300.889
$ws.Range("I105").Value = 2137.5715
$ws.Range("J105").Value = 2872.5
$ws.Range("K105").Value = 2137.5715
$ws.Range("L105").Value = 2872.5
$ws.Range("M105").Value = -390.5715
$ws.Range("N105").Value = -6366.5

$ws.Range("H113").Value = 12346947
$ws.Range("I113").Value = 22223344
$ws.Range("J113").Value = 1449.75
$ws.Range("K113").Value = 22223344
$ws.Range("L113").Value = 1449.75
$ws.Range("M113").Value = -22221174
$ws.Range("N113").Value = -5789.75

$ws.Range("H126").Value = 11115566
$ws.Range("J126").Value = 5483.3335
$ws.Range("L126").Value = 16450.0005
$ws.Range("N126").Value = -21390.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 581452.2
$ws.Range("J5").Value = 703732.2
$ws.Range("L5").Value = 2111196.6
$ws.Range("N5").Value = -2111420.6

$ws.Range("H7").Value = 300
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -1724

$ws.Range("H68").Value = 34066.332
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 50599.5
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 151798.5
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -153420.5

$ws.Range("H71").Value = 34066.332
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 50599.5
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 455395.5
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -463507.5

$ws.Range("H80").Value = 16295.846
$ws.Range("J80").Value = 16828.834
$ws.Range("L80").Value = 50486.50199999999
$ws.Range("N80").Value = -52358.50199999999

$ws.Range("H83").Value = 16295.846
$ws.Range("J83").Value = 16828.834
$ws.Range("L83").Value = 151459.506
$ws.Range("N83").Value = -160819.506

$ws.Range("H92").Value = 374.33334
$ws.Range("J92").Value = 500
$ws.Range("L92").Value = 1500
$ws.Range("N92").Value = -3996

$ws.Range("H94").Value = 1318
$ws.Range("I94").Value = 897.5
$ws.Range("K94").Value = 2692.5
$ws.Range("M94").Value = -2016.5

$ws.Range("H113").Value = 3677086.8
$ws.Range("J113").Value = 7813100
$ws.Range("L113").Value = 23439300
$ws.Range("N113").Value = -23443640

$ws.Range("H122").Value = 2616.8096
$ws.Range("I122").Value = 713.7857
$ws.Range("J122").Value = 3568.3215
$ws.Range("K122").Value = 6424.071300000001
$ws.Range("L122").Value = 32114.8935
$ws.Range("M122").Value = -3974.071300000001
$ws.Range("N122").Value = -37014.8935

$ws.Range("H132").Value = 2124.6667
$ws.Range("J132").Value = 4388.2856
$ws.Range("L132").Value = 39494.5704
$ws.Range("N132").Value = -44554.5704

$ws.Range("H135").Value = 581452.2
$ws.Range("J135").Value = 703732.2
$ws.Range("L135").Value = 6333589.8
$ws.Range("N135").Value = -6338659.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 32794.445
$ws.Range("J88").Value = 32794.445
$ws.Range("L88").Value = 32794.445
$ws.Range("N88").Value = -33696.445

$ws.Range("H91").Value = 32794.445
$ws.Range("J91").Value = 32794.445
$ws.Range("L91").Value = 32794.445
$ws.Range("N91").Value = -35914.445

$ws.Range("H102").Value = 2355.8948
$ws.Range("I102").Value = 1520.4615
$ws.Range("K102").Value = 1520.4615
$ws.Range("M102").Value = 101.5385000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4377.4116
$ws.Range("I7").Value = 3159.889
$ws.Range("J7").Value = 5747.125
$ws.Range("K7").Value = 3159.889
$ws.Range("L7").Value = 5747.125
$ws.Range("M7").Value = -3047.889
$ws.Range("N7").Value = -5971.125

$ws.Range("H126").Value = 4377.4116
$ws.Range("I126").Value = 3159.889
$ws.Range("J126").Value = 5747.125
$ws.Range("K126").Value = 9479.667000000001
$ws.Range("L126").Value = 17241.375
$ws.Range("M126").Value = -7009.667000000001
$ws.Range("N126").Value = -22181.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6250
$ws.Range("I81").Value = 2750
$ws.Range("J81").Value = 9750
$ws.Range("K81").Value = 5500
$ws.Range("L81").Value = 19500
$ws.Range("M81").Value = -4439
$ws.Range("N81").Value = -21622

$ws.Range("H84").Value = 6250
$ws.Range("I84").Value = 2750
$ws.Range("J84").Value = 9750
$ws.Range("K84").Value = 27500
$ws.Range("L84").Value = 97500
$ws.Range("M84").Value = -22196
$ws.Range("N84").Value = -108108

$ws.Range("H122").Value = 3832.0454
$ws.Range("I122").Value = 2608.3076
$ws.Range("J122").Value = 5599.6665
$ws.Range("K122").Value = 7824.9228
$ws.Range("L122").Value = 16798.9995
$ws.Range("M122").Value = -5374.9228
$ws.Range("N122").Value = -21698.9995
